$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure B:E columns remain text-typed so numeric-looking strings (e.g. "2.250", "0.06870")
# keep their exact formatting instead of being coerced to numbers.

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.088.69'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.04%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.872.13'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.35%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.14%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.03'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.14%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.08%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5141'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.98%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3889'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.49%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08385'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.44%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.114'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.28%  '

# Row 11
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.66'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.33%  '

# Row 12
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.191'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.06%  '

# Row 13
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.57'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.08%  '

# Row 14
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.860.07'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.33%  '

# Row 15
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.283'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.82%  '

# Row 16
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.15%  '

# Row 17
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001109'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.71%  '

# Row 18
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.88'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.34%  '

# Row 19
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06647'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.18%  '

# Row 20
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.67'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.44%  '

# Row 21
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.06%  '

# Row 22
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.026'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.20%  '

# Row 23
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'WrappedBTC'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.111.19'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.09%  '

# Row 24
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.12'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.71%  '

# Row 25
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.250'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.41%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.083.70'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.84%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.476'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -4.58%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '157.99'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.82%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.57'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.63%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '124.89'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.03%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1060'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.94%  '

# Row 32
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.57%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.889'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.39%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.596'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.57%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.636'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.35%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02438'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.58%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06525'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.61%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2186'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.33%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.207'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.46%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6503'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.58%  '

# Row 41
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.84%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.225'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.80%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.32'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.48%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6087'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.45%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.03'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.33%  '

# Row 46
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'WEMIXTOKEN'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.279'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.58%  '

# Row 47
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.672'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.24%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.006'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.98%  '

# Row 49
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.59%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '121.37'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.11%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06870'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.56%  '
